# TC28_Verify_PunchOut_User.xlsx
#
# The author's change removes the "Enable Certificate" JavaScript click
# steps (rows 4-8: CLICK_JS / EnableCertificate_MoreInfo / JS_ID / EleType1,
# a WAIT, CLICK_JS / EnableCertificate_GoTOPage / JS_ID / EleType2, and two
# more WAITs) from the "TC28_Verify_PunchOut_User" sheet, which shifts all
# the subsequent rows up by 5 and leaves that sheet as the active one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC28_Verify_PunchOut_User")

# Make this worksheet the active / selected one (it was "Testdata" before).
$ws.Activate()

# Select and delete entire rows 4-8 (the obsolete certificate-enabling
# click steps), shifting rows 9-18 up into rows 4-13.
$ws.Range("A4:A8").EntireRow.Select()
$ws.Range("A4:A8").EntireRow.Delete()
